$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would otherwise be auto-converted to numbers by Excel
# (e.g. "353.85"); force text format first so they stay text like the rest of
# column D, then clear the temporary formatting once the value is set so the
# cell ends up with no explicit style, matching its original unstyled state.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.119.55"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.761.87"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "353.85"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "107.93"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "39.40"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "7.52"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "3.197.80"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "2.760.87"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "0.929"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "51.108.57"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("D20").Value = "3.06"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "69.66"
$ws.Range("D24").Value = "265.18"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "25.94"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").Value = "  +12.01%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "34.54"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "0.0443"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "18.24"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "120.25"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "22.12"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "2.083.51"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").Value = "0.917"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "1.29"
$ws.Range("E51").Value = "  +6.13%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
